$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1424634105014064
$ws.Range("D2").Value = 0.02175719664786158
$ws.Range("E2").Value = 0.4241215750028857
$ws.Range("F2").Value = 0.3861113732231232
$ws.Range("G2").Value = 0.236492415999443
$ws.Range("H2").Value = 0.4088094984408457
$ws.Range("K2").Value = 0.9741754283509181
$ws.Range("O2").Value = 1.205962810138502

$ws.Range("B3").Value = 0.1329242834005129
$ws.Range("D3").Value = 0.01912967109105068
$ws.Range("E3").Value = 0.3700076705874977
$ws.Range("F3").Value = 0.3853025164923878
$ws.Range("G3").Value = 0.2367645389949615
$ws.Range("H3").Value = 0.4131297367877096
$ws.Range("K3").Value = 0.8558119575868375
$ws.Range("O3").Value = 1.215309554227503

$ws.Range("B4").Value = 0.1271362855327709
$ws.Range("D4").Value = 0.0175085037053222
$ws.Range("E4").Value = 0.336863706024829
$ws.Range("F4").Value = 0.3851829900453012
$ws.Range("G4").Value = 0.2372344590999091
$ws.Range("H4").Value = 0.4160588099096785
$ws.Range("K4").Value = 0.7828053301500688
$ws.Range("O4").Value = 1.222259887099895

$ws.Range("B5").Value = 0.1247951600141448
$ws.Range("D5").Value = 0.01684593924354516
$ws.Range("E5").Value = 0.3233763456690184
$ws.Range("F5").Value = 0.3852288202414158
$ws.Range("G5").Value = 0.237501705263746
$ws.Range("H5").Value = 0.4173218430220516
$ws.Range("K5").Value = 0.7529732642226179
$ws.Range("O5").Value = 1.225395757641977

$ws.Range("B6").Value = 0.1244074815375029
$ws.Range("D6").Value = 0.01673580609328695
$ws.Range("E6").Value = 0.3211378853980591
$ws.Range("F6").Value = 0.3852421316731736
$ws.Range("G6").Value = 0.2375506446524938
$ws.Range("H6").Value = 0.4175357581471957
$ws.Range("K6").Value = 0.7480148216396287
$ws.Range("O6").Value = 1.225934769339204

$ws.Range("B7").Value = 0.1271046410967926
$ws.Range("D7").Value = 0.0174995758599934
$ws.Range("E7").Value = 0.3366817356259304
$ws.Range("F7").Value = 0.3851832257325682
$ws.Range("G7").Value = 0.2372377570992015
$ws.Range("H7").Value = 0.4160755627017565
$ws.Range("K7").Value = 0.7824033309388199
$ws.Range("O7").Value = 1.222300950884545

$ws.Range("B8").Value = 0.1391600799395718
$ws.Range("D8").Value = 0.0208528854001031
$ws.Range("E8").Value = 0.4054450657766893
$ws.Range("F8").Value = 0.3857540330290945
$ws.Range("G8").Value = 0.2365231488721236
$ws.Range("H8").Value = 0.4102416844846317
$ws.Range("K8").Value = 0.9334335299433576
$ws.Range("O8").Value = 1.208933502511442

$ws.Range("B9").Value = 0.1633428904134462
$ws.Range("D9").Value = 0.02736450421971881
$ws.Range("E9").Value = 0.5410200769768494
$ws.Range("F9").Value = 0.3898794821020104
$ws.Range("G9").Value = 0.2375430999699759
$ws.Range("H9").Value = 0.4009993260771196
$ws.Range("K9").Value = 1.226905705761055
$ws.Range("O9").Value = 1.192380669631476

$ws.Range("B10").Value = 0.1814345697966928
$ws.Range("D10").Value = 0.03210738695673854
$ws.Range("E10").Value = 0.6411895872659841
$ws.Range("F10").Value = 0.3947633310827214
$ws.Range("G10").Value = 0.2397945820185043
$ws.Range("H10").Value = 0.395555208566833
$ws.Range("K10").Value = 1.440800471857813
$ws.Range("O10").Value = 1.186177663878368

$ws.Range("B11").Value = 0.1897342840725003
$ws.Range("D11").Value = 0.03425569700944209
$ws.Range("E11").Value = 0.6869073979961939
$ws.Range("F11").Value = 0.3973918773714331
$ws.Range("G11").Value = 0.2411506355037147
$ws.Range("H11").Value = 0.3933722529833403
$ws.Range("K11").Value = 1.537719158772234
$ws.Range("O11").Value = 1.184664478650319

$ws.Range("B12").Value = 0.1928870564502319
$ws.Range("D12").Value = 0.03506783243911116
$ws.Range("E12").Value = 0.7042430636519441
$ws.Range("F12").Value = 0.3984460800137768
$ws.Range("G12").Value = 0.2417123002063448
$ws.Range("H12").Value = 0.3925879590608474
$ws.Range("K12").Value = 1.574363020112003
$ws.Range("O12").Value = 1.184280837573198

$ws.Range("B13").Value = 0.1922076151761587
$ws.Range("D13").Value = 0.03489298688430154
$ws.Range("E13").Value = 0.7005084529619836
$ws.Range("F13").Value = 0.3982164170253881
$ws.Range("G13").Value = 0.241589186910673
$ws.Range("H13").Value = 0.3927549855421546
$ws.Range("K13").Value = 1.566473682761739
$ws.Range("O13").Value = 1.184355020259375

$ws.Range("B14").Value = 0.1899934681305382
$ws.Range("D14").Value = 0.03432253989650746
$ws.Range("E14").Value = 0.6883331362822531
$ws.Range("F14").Value = 0.3974774263825793
$ws.Range("G14").Value = 0.2411958763578781
$ws.Range("H14").Value = 0.3933068792990611
$ws.Range("K14").Value = 1.540735028179881
$ws.Range("O14").Value = 1.184629114179387

$ws.Range("B15").Value = 0.1886385168023992
$ws.Range("D15").Value = 0.03397294298643772
$ws.Range("E15").Value = 0.6808784920449824
$ws.Range("F15").Value = 0.3970324441610202
$ws.Range("G15").Value = 0.2409612467023123
$ws.Range("H15").Value = 0.3936504482994252
$ws.Range("K15").Value = 1.524961856987034
$ws.Range("O15").Value = 1.184821701565227

$ws.Range("B16").Value = 0.1808935519251946
$ws.Range("D16").Value = 0.03196679917454048
$ws.Range("E16").Value = 0.6382049902682212
$ws.Range("F16").Value = 0.3945997641036527
$ws.Range("G16").Value = 0.2397126757198436
$ws.Range("H16").Value = 0.3957037884846102
$ws.Range("K16").Value = 1.434458722720763
$ws.Range("O16").Value = 1.186302988507236

$ws.Range("B17").Value = 0.1761599990511229
$ws.Range("D17").Value = 0.03073368789259945
$ws.Range("E17").Value = 0.6120658707757229
$ws.Range("F17").Value = 0.3932118303472762
$ws.Range("G17").Value = 0.2390320226465406
$ws.Range("H17").Value = 0.3970387328906853
$ws.Range("K17").Value = 1.378838390673252
$ws.Range("O17").Value = 1.187547722684116

$ws.Range("B18").Value = 0.1734439573151576
$ws.Range("D18").Value = 0.03002356726354094
$ws.Range("E18").Value = 0.5970453739821835
$ws.Range("F18").Value = 0.3924518037552147
$ws.Range("G18").Value = 0.2386717374743057
$ws.Range("H18").Value = 0.3978341852248874
$ws.Range("K18").Value = 1.346811107432359
$ws.Range("O18").Value = 1.188386760391381

$ws.Range("B19").Value = 0.1725254864420265
$ws.Range("D19").Value = 0.02978298521406231
$ws.Range("E19").Value = 0.5919620497830209
$ws.Range("F19").Value = 0.3922010359199462
$ws.Range("G19").Value = 0.2385550978377466
$ws.Range("H19").Value = 0.3981082530021709
$ws.Range("K19").Value = 1.335961109833988
$ws.Range("O19").Value = 1.188691950193032

$ws.Range("B20").Value = 0.1766632148547131
$ws.Range("D20").Value = 0.03086504487392006
$ws.Range("E20").Value = 0.614846959838232
$ws.Range("F20").Value = 0.3933556144760004
$ws.Range("G20").Value = 0.2391012461925328
$ws.Range("H20").Value = 0.3968937655216038
$ws.Range("K20").Value = 1.384763004481158
$ws.Range("O20").Value = 1.187402469600983

$ws.Range("B21").Value = 0.1906435508353752
$ws.Range("D21").Value = 0.03449013195525197
$ws.Range("E21").Value = 0.691908675127479
$ws.Range("F21").Value = 0.3976928866660785
$ws.Range("G21").Value = 0.2413100907537427
$ws.Range("H21").Value = 0.3931436245268998
$ws.Range("K21").Value = 1.548296661999984
$ws.Range("O21").Value = 1.184543457196781

$ws.Range("B22").Value = 0.1998378450672789
$ws.Range("D22").Value = 0.03685125081102569
$ws.Range("E22").Value = 0.7424097161701013
$ws.Range("F22").Value = 0.4008705960085024
$ws.Range("G22").Value = 0.2430345953459323
$ws.Range("H22").Value = 0.3909395617893665
$ws.Range("K22").Value = 1.65484135671386
$ws.Range("O22").Value = 1.18377927128904

$ws.Range("B23").Value = 0.1949254886011857
$ws.Range("D23").Value = 0.03559183432398072
$ws.Range("E23").Value = 0.7154432731406075
$ws.Range("F23").Value = 0.3991430968496346
$ws.Range("G23").Value = 0.2420883480643568
$ws.Range("H23").Value = 0.3920932815773028
$ws.Range("K23").Value = 1.598007681469653
$ws.Range("O23").Value = 1.184085689269722

$ws.Range("B24").Value = 0.1764356943413503
$ws.Range("D24").Value = 0.03080566208124225
$ws.Range("E24").Value = 0.6135896071114644
$ws.Range("F24").Value = 0.3932904915911024
$ws.Range("G24").Value = 0.2390698535982239
$ws.Range("H24").Value = 0.3969592181348744
$ws.Range("K24").Value = 1.382084643462463
$ws.Range("O24").Value = 1.187467754143768

$ws.Range("B25").Value = 0.1567433744421862
$ws.Range("D25").Value = 0.0256100274784572
$ws.Range("E25").Value = 0.5042531803391768
$ws.Range("F25").Value = 0.3884395160511573
$ws.Range("G25").Value = 0.2370052668188265
$ws.Range("H25").Value = 0.4032636759672812
$ws.Range("K25").Value = 1.147810306434053
$ws.Range("O25").Value = 1.195817278744627
